$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.586.47'
$ws.Range("E2").Value = '  -0.05%  '

$ws.Range("D3").Value = '1.936.97'
$ws.Range("E3").Value = '  +0.70%  '

$ws.Range("D4").Value = "'0.9989"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").Value = "'246.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.58%  '

$ws.Range("D6").Value = "'0.9995"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.09%  '

$ws.Range("D7").Value = "'0.4838"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.29%  '

$ws.Range("D8").Value = "'0.2926"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.20%  '

$ws.Range("D9").Value = "'0.06821"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.28%  '

$ws.Range("D10").Value = "'113.52"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.28%  '

$ws.Range("D11").Value = "'19.46"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.45%  '

$ws.Range("D12").Value = '1.936.21'
$ws.Range("E12").Value = '  +0.64%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = "'5.499"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.04%  '

$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").Value = "'0.07594"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.76%  '

$ws.Range("D15").Value = "'0.6825"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.36%  '

$ws.Range("D16").Value = "'298.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.60%  '

$ws.Range("D17").Value = '30.588.78'

$ws.Range("E18").Value = '  +1.31%  '

$ws.Range("D19").Value = "'0.000007679"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.36%  '

$ws.Range("D20").Value = "'5.580"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.02%  '

$ws.Range("D21").Value = "'0.9988"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.17%  '

$ws.Range("D22").Value = '2.187.64'
$ws.Range("E22").Value = '  +0.47%  '

$ws.Range("D23").Value = "'0.9991"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.12%  '

$ws.Range("D24").Value = "'6.528"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.43%  '

$ws.Range("D25").Value = "'9.555"
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").Value = "'168.03"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.44%  '

$ws.Range("D27").Value = "'20.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.78%  '

$ws.Range("D28").Value = "'2.142"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.10%  '

$ws.Range("E29").Value = '  -0.06%  '

$ws.Range("D30").Value = "'1.434"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.86%  '

$ws.Range("E31").Value = '  -0.83%  '

$ws.Range("D32").Value = "'4.113"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.12%  '

$ws.Range("D33").Value = "'0.05016"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.70%  '

$ws.Range("D34").Value = "'0.7481"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.63%  '

$ws.Range("E35").Value = '  +0.41%  '

$ws.Range("D36").Value = "'0.02042"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.08%  '

$ws.Range("D37").Value = "'2.713"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.30%  '

$ws.Range("D38").Value = "'2.699"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.52%  '

$ws.Range("D39").Value = "'2.044"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.40%  '

$ws.Range("D40").Value = "'110.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.44%  '

$ws.Range("D41").Value = "'0.4467"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.28%  '

$ws.Range("E42").Value = '  +0.14%  '

$ws.Range("D43").Value = "'5.863"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.87%  '

$ws.Range("D44").Value = "'69.93"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.78%  '

$ws.Range("E45").Value = '  +0.07%  '

$ws.Range("D46").Value = "'7.311"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.17%  '

$ws.Range("D47").Value = "'49.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.07%  '

$ws.Range("D48").Value = "'9.360"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.28%  '

$ws.Range("D49").Value = "'0.1236"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.99%  '

$ws.Range("D50").Value = "'0.2553"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.79%  '

$ws.Range("D51").Value = "'35.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.66%  '
